# Update "想去人数" (number of people interested) figures for three
# exhibition entries that appear on both the "展览" sheet and the
# combined "全部类型" sheet.
#
#   Row 4 (布谷鸟动漫展5th):            120  -> 124
#   Row 5 (2024良牙动漫秋季盛典（秋典）): 2779 -> 2798
#   Row 6 (快看漫画动漫游戏嘉年华):       271  -> 275

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 124
    $ws.Range("F5").Value = 2798
    $ws.Range("F6").Value = 275
}
